# process executive data in test2.py
# - E-1 (Target/Sales) figures for Brand4 (row 5) are filled in: B5=20, C5=25.
# - The previously-entered E-2 figures for Brand6 (row 8, D8/E8) and the
#   E-3 figures for Brand7 (row 9, F9/G9) are cleared back to 0, and their
#   cells revert to the plain (non-bold, default-height) row style.
# - Selection moves to the newly-relevant range D8:E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Brand4): fill E-1 Target / E-1 Sales, restyle to the bold xf
#     that D5/E5 already use (font size 12) instead of the near-duplicate
#     bold-11 font so the workbook's font table collapses back to 2 fonts.
$ws.Range("B5:C5").Font.Size = 12
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 25

# --- Row 8 (Brand6): clear E-2 Target / E-2 Sales back to 0 and drop the
#     bold styling so the cells fall back to the plain bordered xf.
$ws.Range("D8:E8").Font.Size = 11
$ws.Range("D8:E8").Font.Bold = $false
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
[void]$ws.Rows.Item(8).AutoFit()

# --- Row 9 (Brand7): clear E-3 Target / E-3 Sales back to 0 and drop the
#     bold styling so the cells fall back to the plain bordered xf.
$ws.Range("F9:G9").Font.Size = 11
$ws.Range("F9:G9").Font.Bold = $false
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
[void]$ws.Rows.Item(9).AutoFit()

# --- Move the active selection to the range touched last (D8:E8).
[void]$ws.Range("D8:E8").Select()
